# edit.ps1 -- applies the "Celestial Symphony" -> "The Enigma of Government"
# content rewrite described by the target diff, via Word COM interop.

$d = $word.ActiveDocument
$br = [char]11   # Word "manual line break" -> serializes as <w:br/>

## --- Title ---
$r = $d.Content
[void]$r.Find.Execute("Celestial Symphony: The Cosmic Dance of Harmony", $true, $false, $false, $false, $false, $true, 1, $false, "The Enigma of Government: Unraveling the Complex Web of Governance", 2)

## --- Author name ---
$r = $d.Content
[void]$r.Find.Execute("[Valid Author's Name]", $true, $false, $false, $false, $false, $true, 1, $false, "Alex Carter", 2)

## --- Email line (scoped to its paragraph to avoid stray matches) ---
$r = $d.Paragraphs(3).Range
[void]$r.Find.Execute("Valid Email Address - e", $true, $false, $false, $false, $false, $true, 1, $false, "alex", 2)

$r = $d.Paragraphs(3).Range
[void]$r.Find.Execute("g.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$r = $d.Paragraphs(3).Range
[void]$r.Find.Execute(" username@domainname", $true, $false, $false, $false, $false, $true, 1, $false, "carter938@educators", 2)

## --- Body paragraph: opening block ---
$r = $d.Content
[void]$r.Find.Execute("In the vast expanse of the cosmos, amidst the tapestry of celestial bodies, lies a mesmerizing interplay of rhythm and harmony, an intricate cosmic dance that mirrors the very essence of life on Earth", $true, $false, $false, $false, $false, $true, 1, $false, "Government, an intricate tapestry of intricate mechanisms and processes, stands as a bastion of societal organization", 2)

$r = $d.Content
[void]$r.Find.Execute(" From the gravitational ballet of planets to the synchronized pulsars, each celestial entity contributes its unique voice to a symphony that resounds across the universe", $true, $false, $false, $false, $false, $true, 1, $false, " It is a labyrinthine network of institutions, laws, and policies designed to regulate human conduct, facilitate social harmony, and secure the general welfare", 2)

$r = $d.Content
[void]$r.Find.Execute(" It is a labyrinthine network of institutions, laws, and policies designed to regulate human conduct, facilitate social harmony, and secure the general welfare", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(". Throughout history, governments have taken various forms, ranging from monarchies to democracies, each with its unique characteristics and challenges. Understanding the enigmatic nature of governance is vital for responsible citizenship in a democratic society. This essay delves into the enigma of government, deciphering the intricate mechanisms that shape our political landscape")

## --- Body paragraph: second block ---
$r = $d.Content
[void]$r.Find.Execute("From the celestial bodies orbiting the Sun in perfect harmony like musicians following a cosmic score, their gravitational dance shapes the solar system's dynamics", $true, $false, $false, $false, $false, $true, 1, $false, "Governments, acting as the custodians of public interest, are entrusted with the formidable task of managing societal affairs", 2)

$r = $d.Content
[void]$r.Find.Execute(" Their periodic movements, governed by celestial mechanics, create a rhythmic flow, a delicate balance that guides the symphony of the spheres", $true, $false, $false, $false, $false, $true, 1, $false, " They are responsible for formulating and implementing policies that touch every aspect of our lives, from taxation to education to national defense", 2)

$r = $d.Content
[void]$r.Find.Execute(" Similarly, in the realm of subatomic particles, the harmonious vibrations of atoms and molecules set the stage for the intricate dance of matter, shaping the physical world we experience", $true, $false, $false, $false, $false, $true, 1, $false, " Balancing the needs of diverse constituencies with the limitations of available resources is a perpetual challenge that governments face", 2)

$r = $d.Content
[void]$r.Find.Execute(" Balancing the needs of diverse constituencies with the limitations of available resources is a perpetual challenge that governments face", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(". The allocation of resources, regulation of industries, and provision of public services are just a few of the responsibilities that lie within the government's purview")

## --- Body paragraph: third block + new "Body:" section ---
$r = $d.Content
[void]$r.Find.Execute("Beyond our solar system, stars engage in a captivating cosmic ballet, forming celestial clusters that shine in unison", $true, $false, $false, $false, $false, $true, 1, $false, "Furthermore, governments play a crucial role in maintaining social order and upholding the rule of law", 2)

$r = $d.Content
[void]$r.Find.Execute(" Their synchronized pulsations, emitting regular bursts of energy, create an intricate rhythm in the cosmic tapestry, akin to the harmonious chords of a celestial symphony", $true, $false, $false, $false, $false, $true, 1, $false, " Through its law enforcement agencies, judicial system, and regulatory bodies, a government ensures that citizens' rights and safety are protected", 2)

$r = $d.Content
[void]$r.Find.Execute(" These stellar collaborations showcase the delicate balance and order amidst the vastness of space", $true, $false, $false, $false, $false, $true, 1, $false, " It establishes a framework for resolving disputes, enforcing contracts, and deterring criminal behavior", 2)

$r = $d.Content
[void]$r.Find.Execute(" It establishes a framework for resolving disputes, enforcing contracts, and deterring criminal behavior", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(". By doing so, governments foster an environment conducive to economic prosperity, social progress, and personal freedom." + $br + $br + "Body:" + $br + $br + "History bears witness to the ever-evolving nature of governance, reflecting humanity's quest for more just and equitable societies. From the ancient city-states of Greece to the sprawling empires of Rome and China, governments have undergone profound transformations, shaped by cultural, economic, and political forces. These changes have given rise to a multitude of governance models, each with its strengths and weaknesses. Democracy, characterized by the participation of citizens in decision-making, has emerged as a popular choice, although it remains vulnerable to challenges such as populism and interest group influence." + $br + $br + "In contemporary times, governments are grappling with a multitude of challenges, including globalization, technological advancements, and climate change. The interconnectedness of the global economy has created a complex web of interdependence, requiring governments to collaborate and coordinate policies on issues such as trade, finance, and environmental sustainability. Technological advancements have introduced new ethical and regulatory dilemmas, demanding proactive responses from governments. Climate change, a pressing global crisis, poses unprecedented challenges to governance, requiring international cooperation and concerted action to mitigate its devastating effects." + $br + $br + "Despite these challenges, governments remain the cornerstone of organized societies, providing essential services, upholding the rule of law, and facilitating economic progress. By understanding the enigmatic nature of governance, citizens can engage more meaningfully in the political process, hold their leaders accountable, and work towards building a better future for themselves and succeeding generations")

## --- Summary paragraph ---
$r = $d.Content
[void]$r.Find.Execute("The cosmic dance of harmony extends beyond our solar system, with distant galaxies and stellar clusters playing their own musical compositions in the grand symphony of the universe", $true, $false, $false, $false, $false, $true, 1, $false, "Governments, as complex and enigmatic entities, play a vital role in managing societal affairs, maintaining social order, and securing the general welfare", 2)

$r = $d.Content
[void]$r.Find.Execute(" From the synchronized pulsations of pulsars to the graceful waltz of orbiting planets, the cosmos is a symphony of movements, rhythms, and harmonies that mirrors the complexities of life on Earth", $true, $false, $false, $false, $false, $true, 1, $false, " Throughout history, they have undergone profound transformations, adapting to changing cultural, economic, and political landscapes", 2)

$r = $d.Content
[void]$r.Find.Execute(" The underlying interconnectedness of everything invites us to contemplate the unity between the seemingly disparate realms of science, art, music, and philosophy, revealing the universal language of harmony that resonates within and beyond our world", $true, $false, $false, $false, $false, $true, 1, $false, " Today, governments face a multitude of challenges, including globalization, technological advancements, and climate change", 2)

$r = $d.Content
[void]$r.Find.Execute(" Today, governments face a multitude of challenges, including globalization, technological advancements, and climate change", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(". Understanding the enigma of government empowers citizens to participate more effectively in the political process and work towards building a more just and prosperous society")

## --- Trailing empty paragraph ---
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

